$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 150
$ws.Range("I33").Value = 150
$ws.Range("K33").Value = 150
$ws.Range("M33").Value = 79
$ws.Range("H51").Value = 8420
$ws.Range("J51").Value = 8420
$ws.Range("L51").Value = 8420
$ws.Range("N51").Value = -9388
$ws.Range("H103").Value = 1700
$ws.Range("I103").Value = 1700
$ws.Range("K103").Value = 5100
$ws.Range("M103").Value = -4514

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 80000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 80000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 80000
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -82372
$ws.Range("H89").Value = 80000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 80000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 240000
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -251856
$ws.Range("H122").Value = 1999.7142
$ws.Range("I122").Value = 1808.7273
$ws.Range("K122").Value = 5426.1819
$ws.Range("M122").Value = -2976.1819
$ws.Range("H132").Value = 2173.0833
$ws.Range("I132").Value = 1807.7
$ws.Range("K132").Value = 5423.1
$ws.Range("M132").Value = -2893.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 390.5
$ws.Range("I64").Value = 390.5
$ws.Range("K64").Value = 390.5
$ws.Range("M64").Value = -165.5
$ws.Range("H67").Value = 390.5
$ws.Range("I67").Value = 390.5
$ws.Range("K67").Value = 390.5
$ws.Range("M67").Value = 389.5
$ws.Range("H94").Value = 4698.778
$ws.Range("I94").Value = 3414.8333
$ws.Range("J94").Value = 7266.6665
$ws.Range("K94").Value = 3414.8333
$ws.Range("L94").Value = 7266.6665
$ws.Range("M94").Value = -2963.8333
$ws.Range("N94").Value = -8168.6665
$ws.Range("H105").Value = 4934.8335
$ws.Range("I105").Value = 5005
$ws.Range("J105").Value = 4899.75
$ws.Range("K105").Value = 5005
$ws.Range("L105").Value = 4899.75
$ws.Range("M105").Value = -3258
$ws.Range("N105").Value = -8393.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 4
$ws.Range("I19").Value = 4
$ws.Range("K19").Value = 4
$ws.Range("M19").Value = 166
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = 4
$ws.Range("K24").Value = 4
$ws.Range("M24").Value = 166
$ws.Range("H31").Value = 4488.4614
$ws.Range("I31").Value = 1925.8572
$ws.Range("K31").Value = 1925.8572
$ws.Range("M31").Value = -1630.8572
$ws.Range("H34").Value = 4488.4614
$ws.Range("I34").Value = 1925.8572
$ws.Range("K34").Value = 1925.8572
$ws.Range("M34").Value = -1723.8572
$ws.Range("H58").Value = 4494.75
$ws.Range("I58").Value = 4494.75
$ws.Range("K58").Value = 4494.75
$ws.Range("M58").Value = -4291.75
$ws.Range("H99").Value = 3900
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 3900
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 3900
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = -6896
$ws.Range("H106").Value = 69998
$ws.Range("J106").Value = 69998
$ws.Range("L106").Value = 69998
$ws.Range("N106").Value = -72522
$ws.Range("H126").Value = 3900
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3900
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 11700
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -16640
$ws.Range("H136").Value = 4494.75
$ws.Range("I136").Value = 4494.75
$ws.Range("K136").Value = 13484.25
$ws.Range("M136").Value = -10934.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 1550
$ws.Range("J123").Value = 2500
$ws.Range("L123").Value = 7500
$ws.Range("N123").Value = -12400

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 566100.25
$ws.Range("J3").Value = 1007500
$ws.Range("L3").Value = 1007500
$ws.Range("N3").Value = -1007732
$ws.Range("H31").Value = 523.25
$ws.Range("I31").Value = 523.25
$ws.Range("K31").Value = 523.25
$ws.Range("M31").Value = -231.25
$ws.Range("H37").Value = 523.25
$ws.Range("I37").Value = 523.25
$ws.Range("K37").Value = 523.25
$ws.Range("M37").Value = -246.25
$ws.Range("H80").Value = 2337.375
$ws.Range("J80").Value = 2499.6667
$ws.Range("L80").Value = 2499.6667
$ws.Range("N80").Value = -4495.6667
$ws.Range("H83").Value = 2337.375
$ws.Range("J83").Value = 2499.6667
$ws.Range("L83").Value = 12498.3335
$ws.Range("N83").Value = -22482.3335
$ws.Range("H102").Value = 2238.5
$ws.Range("I102").Value = 2289.0557
$ws.Range("J102").Value = 1783.5
$ws.Range("K102").Value = 2289.0557
$ws.Range("L102").Value = 1783.5
$ws.Range("M102").Value = -667.0556999999999
$ws.Range("N102").Value = -5027.5
$ws.Range("H122").Value = 1524
$ws.Range("I122").Value = 1661.3334
$ws.Range("K122").Value = 4984.0002
$ws.Range("M122").Value = -2534.0002
$ws.Range("H126").Value = 1706
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 9000
$ws.Range("I25").Value = 9000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 9000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -8770
$ws.Range("N25").Value = ""
$ws.Range("H46").Value = 4323.205
$ws.Range("I46").Value = 2565.6667
$ws.Range("K46").Value = 2565.6667
$ws.Range("M46").Value = -2377.6667
$ws.Range("H122").Value = 7298.3335
$ws.Range("I122").Value = 7848.5
$ws.Range("K122").Value = 23545.5
$ws.Range("M122").Value = -21095.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5105000.5
$ws.Range("I2").Value = 5105000.5
$ws.Range("K2").Value = 5105000.5
$ws.Range("M2").Value = -5104888.5
$ws.Range("H4").Value = 45969.445
$ws.Range("I4").Value = 67638.664
$ws.Range("K4").Value = 67638.664
$ws.Range("M4").Value = -67525.664
$ws.Range("H5").Value = 6666800
$ws.Range("I5").Value = 401
$ws.Range("K5").Value = 401
$ws.Range("M5").Value = -289
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = ""
$ws.Range("H122").Value = 2176.2
$ws.Range("I122").Value = 2176.2
$ws.Range("K122").Value = 6528.599999999999
$ws.Range("M122").Value = -4078.599999999999
$ws.Range("H126").Value = 5049
$ws.Range("J126").Value = 1697
$ws.Range("L126").Value = 5091
$ws.Range("N126").Value = -10031
